# Results from October_25,_2020--23:17:44 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateCell($addr, $dateSerial) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "YYYY-MM-DD"
    $c.Value = $dateSerial
}

function Set-TextDateCell($addr, $text) {
    # Forces the literal string to be stored as text (not re-interpreted
    # as an Excel date serial number), while leaving the cell with the
    # default (no explicit) style.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------
# Row 5 - North Carolina
# ---------------------------------------------------------------------
Set-DateCell "B5" 44129
$ws.Range("C5").Value = 260099
$ws.Range("D5").Value = 4157
$ws.Range("E5").Value = 47548
$ws.Range("F5").Value = 1184
$ws.Range("G5").Value = 22.74
$ws.Range("H5").Value = 29.94
$ws.Range("J5").Value = $true
$ws.Range("K5").Value = 209101
$ws.Range("L5").Value = 3954
$ws.Range("O5").Value = "Success!"

# ---------------------------------------------------------------------
# Row 6 - Wyoming
# ---------------------------------------------------------------------
$ws.Range("O6").Value = "An error occurred. ... KeyError('SUM(# probable)')"

# ---------------------------------------------------------------------
# Row 7 - NewYork
# ---------------------------------------------------------------------
$ws.Range("O7").Value = "An error occurred. ... TimeoutException('', None, None)"

# ---------------------------------------------------------------------
# Row 11 - Oregon
# ---------------------------------------------------------------------
$ws.Range("O11").Value = 'An error occurred. ... KeyError("None of [''Categories''] are in the columns")'

# ---------------------------------------------------------------------
# Row 13 - SouthCarolina
# ---------------------------------------------------------------------
$ws.Range("O13").Value = 'An error occurred. ... AttributeError("''NoneType'' object has no attribute ''text''")'

# ---------------------------------------------------------------------
# Row 16 - Maryland
# ---------------------------------------------------------------------
Set-DateCell "B16" 44129
$ws.Range("C16").Value = 140279
$ws.Range("D16").Value = 3950
$ws.Range("E16").Value = 43513
$ws.Range("F16").Value = 1606
$ws.Range("G16").Value = 36.45
$ws.Range("H16").Value = 40.78
$ws.Range("K16").Value = 119377
$ws.Range("L16").Value = 3938
$ws.Range("O16").Value = "Success!"

# ---------------------------------------------------------------------
# Row 19 - Ohio
# ---------------------------------------------------------------------
$ws.Range("O19").Value = "An error occurred. ... JSONDecodeError('Expecting value: line 1 column 1 (char 0)')"

# ---------------------------------------------------------------------
# Row 23 - Oklahoma (note: B23 stays literal text, not a real date)
# ---------------------------------------------------------------------
Set-TextDateCell "B23" "2020-10-25"
$ws.Range("C23").Value = 116736
$ws.Range("D23").Value = 1249
$ws.Range("E23").Value = 8031.4368
$ws.Range("F23").Value = 80.0609
$ws.Range("G23").Value = 6.88
$ws.Range("H23").Value = 6.41
$ws.Range("K23").Value = 95151.51360000001
$ws.Range("L23").Value = 1157.9479
$ws.Range("O23").Value = "Success!"

# ---------------------------------------------------------------------
# Row 25 - Kansas
# ---------------------------------------------------------------------
$ws.Range("O25").Value = "An error occurred. ... TimeoutException('', None, None)"

# ---------------------------------------------------------------------
# Row 31 - Iowa
# ---------------------------------------------------------------------
Set-DateCell "B31" 44130
$ws.Range("C31").Value = 116238
$ws.Range("D31").Value = 1635
$ws.Range("E31").Value = 5488
$ws.Range("F31").Value = 57
$ws.Range("G31").Value = 4.72
$ws.Range("H31").Value = 3.49
$ws.Range("I31").Value = $true
$ws.Range("J31").Value = $true
$ws.Range("O31").Value = "Success!"

# ---------------------------------------------------------------------
# Row 37 - Nevada
# ---------------------------------------------------------------------
$ws.Range("O37").Value = "An error occurred. ... TimeoutException('', None, None)"

# ---------------------------------------------------------------------
# Row 39 - Delaware
# ---------------------------------------------------------------------
$ws.Range("O39").Value = 'An error occurred. ... NoSuchElementException(''no such element: Unable to locate element: {"method":"xpath","selector":"//a[@data-chart-id=\''count-charts\'']"}\n  (Session info: headless chrome=86.0.4240.111)'', None, None)'

# ---------------------------------------------------------------------
# Row 42 - SouthDakota
# ---------------------------------------------------------------------
$ws.Range("O42").Value = "An error occurred. ... TimeoutException('', None, None)"

# ---------------------------------------------------------------------
# Row 47 - California - San Francisco
# ---------------------------------------------------------------------
$ws.Range("O47").Value = "An error occurred. ... TimeoutException('', None, None)"

# ---------------------------------------------------------------------
# Row 50 - Idaho
# ---------------------------------------------------------------------
$ws.Range("O50").Value = "An error occurred. ... TimeoutException('', None, None)"

# ---------------------------------------------------------------------
# Row 52 - Arizona
# ---------------------------------------------------------------------
Set-DateCell "B52" 44129
$ws.Range("C52").Value = 238163
$ws.Range("D52").Value = 5874
$ws.Range("E52").Value = 7371
$ws.Range("F52").Value = 181
$ws.Range("G52").Value = 4.39
$ws.Range("H52").Value = 3.46
$ws.Range("K52").Value = 167906
$ws.Range("L52").Value = 5235
$ws.Range("O52").Value = "Success!"

# ---------------------------------------------------------------------
# Row 55 - WestVirginia
# ---------------------------------------------------------------------
$ws.Range("O55").Value = "An error occurred. ... TimeoutException('', None, None)"

# ---------------------------------------------------------------------
# Row 57 - New Hampshire
# ---------------------------------------------------------------------
Set-DateCell "B57" 44129
$ws.Range("C57").Value = 10328
$ws.Range("D57").Value = 473
$ws.Range("E57").Value = 408
$ws.Range("F57").Value = 9
$ws.Range("G57").Value = 4.55
$ws.Range("H57").Value = 1.99
$ws.Range("K57").Value = 8964
$ws.Range("L57").Value = 452
$ws.Range("O57").Value = "Success!"

Write-Host "Edit complete"
